$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - match formatting of existing header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows (2-6) for new columns I and J
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 6
